# Add the new "Asha Enterprise" outstanding entry (row 24) to the
# "Purchase 22-23" sheet, and update the active sheet / selections on
# both sheets to match where the user left off after the edit.

$wb = $excel.ActiveWorkbook

$wsPurchase = $wb.Worksheets.Item("Purchase 22-23")
$wsSale = $wb.Worksheets.Item("Sale 22-23")

# --- Add new row 24 on "Purchase 22-23" sheet ---
# Copy the formatting of the previous data row (row 22) down onto the new
# row first, then fill in the new record's values/formula.
$wsPurchase.Range("A22:F22").Copy()
$wsPurchase.Range("A24:F24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsPurchase.Range("A24").Value = 10
$wsPurchase.Range("B24").Value = 45295
$wsPurchase.Range("C24").Value = 2117
$wsPurchase.Range("D24").Value = "Asha Enterprise"
$wsPurchase.Range("E24").Value = 1569
$wsPurchase.Range("F24").Formula = "=E24"

# --- Selection / active sheet updates ---
# "Sale 22-23" keeps a lingering selection even though it's no longer active.
$wsSale.Range("D27:E27").Select()

# "Purchase 22-23" becomes the active (visible) sheet, selection moved below
# the freshly added row.
$wsPurchase.Activate()
$wsPurchase.Range("A25").Select()
